# Update "want to go" counts (column F) on the "展览" (Exhibition) and
# "演出" (Performance) sheets, as well as the aggregated "全部类型"
# (All Types) sheet which mirrors the same rows.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 126
$wsExhibition.Range("F4").Value = 167
$wsExhibition.Range("F5").Value = 3210
$wsExhibition.Range("F6").Value = 325
$wsExhibition.Range("F8").Value = 416

# --- 演出 (Performance) sheet ---
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F3").Value = 4

# --- 全部类型 (All Types) sheet, mirrors rows from both sheets above ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 126
$wsAll.Range("F4").Value = 167
$wsAll.Range("F5").Value = 3210
$wsAll.Range("F6").Value = 325
$wsAll.Range("F8").Value = 4
$wsAll.Range("F10").Value = 416
